# Updated symbol list on Sun Jan 29 13:40:50 UTC 2023 with GitHub Actions
# Applies refreshed price / volume(1h) figures to the cryptos sheet,
# including a position swap between BOLO and CoinbaseStockToken (rows 48-49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '314.69' },
    @{ Cell = 'E2'; Value = '3.23%' },
    @{ Cell = 'D3'; Value = '39.42' },
    @{ Cell = 'E3'; Value = '2.71%' },
    @{ Cell = 'D4'; Value = '5.113' },
    @{ Cell = 'E4'; Value = '0.49%' },
    @{ Cell = 'D5'; Value = '0.08200' },
    @{ Cell = 'E5'; Value = '1.86%' },
    @{ Cell = 'D6'; Value = '1.966' },
    @{ Cell = 'E6'; Value = '1.95%' },
    @{ Cell = 'E7'; Value = '3.74%' },
    @{ Cell = 'D8'; Value = '0.9306' },
    @{ Cell = 'E8'; Value = '0.05%' },
    @{ Cell = 'D9'; Value = '0.1415' },
    @{ Cell = 'E9'; Value = '-1.34%' },
    @{ Cell = 'D10'; Value = '0.1971' },
    @{ Cell = 'E10'; Value = '1.98%' },
    @{ Cell = 'D11'; Value = '0.09135' },
    @{ Cell = 'E11'; Value = '1.33%' },
    @{ Cell = 'D12'; Value = '0.03539' },
    @{ Cell = 'E12'; Value = '0.36%' },
    @{ Cell = 'D13'; Value = '0.09821' },
    @{ Cell = 'E13'; Value = '0.53%' },
    @{ Cell = 'D14'; Value = '0.001402' },
    @{ Cell = 'E14'; Value = '0.48%' },
    @{ Cell = 'D15'; Value = '0.005945' },
    @{ Cell = 'E15'; Value = '-1.75%' },
    @{ Cell = 'D16'; Value = '3.661' },
    @{ Cell = 'E16'; Value = '-1.73%' },
    @{ Cell = 'D17'; Value = '4.269' },
    @{ Cell = 'E17'; Value = '1.70%' },
    @{ Cell = 'D18'; Value = '3.322' },
    @{ Cell = 'E18'; Value = '-3.03%' },
    @{ Cell = 'D19'; Value = '0.3462' },
    @{ Cell = 'E19'; Value = '0.03%' },
    @{ Cell = 'D20'; Value = '0.1291' },
    @{ Cell = 'E20'; Value = '-3.19%' },
    @{ Cell = 'D21'; Value = '4.893' },
    @{ Cell = 'E21'; Value = '1.09%' },
    @{ Cell = 'D22'; Value = '0.2446' },
    @{ Cell = 'E22'; Value = '1.56%' },
    @{ Cell = 'D23'; Value = '0.04315' },
    @{ Cell = 'E23'; Value = '-1.20%' },
    @{ Cell = 'D24'; Value = '0.001222' },
    @{ Cell = 'E24'; Value = '-0.54%' },
    @{ Cell = 'E25'; Value = '16.39%' },
    @{ Cell = 'D26'; Value = '0.0001298' },
    @{ Cell = 'E26'; Value = '-0.39%' },
    @{ Cell = 'D27'; Value = '0.0003994' },
    @{ Cell = 'E27'; Value = '-10.19%' },
    @{ Cell = 'D39'; Value = '0.02240' },
    @{ Cell = 'E39'; Value = '7.58%' },
    @{ Cell = 'D40'; Value = '0.05277' },
    @{ Cell = 'E40'; Value = '4.79%' },
    @{ Cell = 'D41'; Value = '0.007566' },
    @{ Cell = 'E41'; Value = '1.05%' },
    @{ Cell = 'D42'; Value = '0.009862' },
    @{ Cell = 'E42'; Value = '-2.66%' },
    @{ Cell = 'E43'; Value = '2.59%' },
    @{ Cell = 'D44'; Value = '0.002116' },
    @{ Cell = 'E44'; Value = '-1.39%' },
    @{ Cell = 'D45'; Value = '0.009802' },
    @{ Cell = 'E45'; Value = '10.24%' },
    @{ Cell = 'D46'; Value = '0.00006362' },
    @{ Cell = 'E46'; Value = '2.81%' },
    @{ Cell = 'D47'; Value = '0.00000000749' },
    @{ Cell = 'E47'; Value = '-0.34%' },
    @{ Cell = 'B48'; Value = 'BOLO' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo' },
    @{ Cell = 'D48'; Value = '0.002763' },
    @{ Cell = 'E48'; Value = '-7.50%' },
    @{ Cell = 'B49'; Value = 'CoinbaseStockToken' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin' },
    @{ Cell = 'D49'; Value = '0.001198' },
    @{ Cell = 'E49'; Value = '-25.20%' },
    @{ Cell = 'D50'; Value = '0.00002097' },
    @{ Cell = 'E50'; Value = '-0.34%' },
    @{ Cell = 'D51'; Value = '0.0001997' },
    @{ Cell = 'E51'; Value = '-0.34%' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}

